$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (columns B, C, D, E, F)
$data = @(
    @("NSE:ADROITINFO", "NSE:BEARDSELL", "NSE:CIPLA",      "", "NSE:NESTLEIND"),
    @("NSE:BAJAJCON",   "NSE:DCMSRIND",  "NSE:GLENMARK",   "", ""),
    @("NSE:CAPLIPOINT", "NSE:ESTER",     "NSE:GRASIM",     "", ""),
    @("NSE:COASTCORP",  "NSE:GENCON",    "NSE:LAURUSLABS", "", ""),
    @("NSE:CYIENTDLM",  "NSE:HBSL",      "NSE:LT",         "", ""),
    @("NSE:HCG",        "NSE:INDIANHUME","NSE:M&M",        "", ""),
    @("NSE:KRISHANA",   "NSE:INGERRAND", "NSE:NESTLEIND",  "", ""),
    @("NSE:NAGREEKCAP", "NSE:IPL",       "NSE:PPLPHARMA",  "", ""),
    @("NSE:NESTLEIND",  "NSE:JWL",       "",               "", ""),
    @("NSE:PALREDTEC",  "NSE:KIOCL",     "",               "", ""),
    @("",               "NSE:KIRLPNU",   "",               "", ""),
    @("",               "NSE:NEWGEN",    "",               "", ""),
    @("",               "NSE:RANEHOLDIN","",               "", ""),
    @("",               "NSE:RHIM",      "",               "", ""),
    @("",               "NSE:RUCHIRA",   "",               "", "")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Range("B$r").Value = $rowVals[0]
    $ws.Range("C$r").Value = $rowVals[1]
    $ws.Range("D$r").Value = $rowVals[2]
    $ws.Range("E$r").Value = $rowVals[3]
    $ws.Range("F$r").Value = $rowVals[4]
}

# Remove the now-unused rows 17-23 so the sheet dimension shrinks to A1:F16
$ws.Range("A17:F23").EntireRow.Delete()
